$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.07002266666666
$ws.Range("H2").Value = 36.21006799999999
$ws.Range("I2").Value = 0.7601982364861632
$ws.Range("J2").Value = 0.7601982364861634
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 62.90731233333333
$ws.Range("N2").Value = 188.721937
$ws.Range("O2").Value = 0.9393635410440488
$ws.Range("P2").Value = 0.9393635410440487
$ws.Range("Q2").Value = 759.2926857624127
$ws.Range("R2").Value = 6833.634171861714
$ws.Range("S2").Value = 0.7141025073210836
$ws.Range("T2").Value = 0.7141025073210836

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.07002266666666
$ws.Range("H3").Value = 36.21006799999999
$ws.Range("I3").Value = 0.7601982364861632
$ws.Range("J3").Value = 0.7601982364861634
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.253965666666666
$ws.Range("N3").Value = 6.761896999999999
$ws.Range("O3").Value = 0.03365734588711396
$ws.Range("P3").Value = 0.03365734588711396
$ws.Range("Q3").Value = 27.2054166865551
$ws.Range("R3").Value = 244.8487501789959
$ws.Range("S3").Value = 0.02558625498818885
$ws.Range("T3").Value = 0.02558625498818885

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.07002266666666
$ws.Range("H4").Value = 36.21006799999999
$ws.Range("I4").Value = 0.7601982364861632
$ws.Range("J4").Value = 0.7601982364861634
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.806737666666667
$ws.Range("N4").Value = 5.420213
$ws.Range("O4").Value = 0.02697911306883729
$ws.Range("P4").Value = 0.02697911306883729
$ws.Range("Q4").Value = 21.80736458938711
$ws.Range("R4").Value = 196.266281304484
$ws.Range("S4").Value = 0.02050947417689091
$ws.Range("T4").Value = 0.02050947417689091

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.308268
$ws.Range("H5").Value = 3.924804
$ws.Range("I5").Value = 0.08239777620284613
$ws.Range("J5").Value = 0.08239777620284613
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 62.90731233333333
$ws.Range("N5").Value = 188.721937
$ws.Range("O5").Value = 0.9393635410440488
$ws.Range("P5").Value = 0.9393635410440487
$ws.Range("Q5").Value = 82.29962369170532
$ws.Range("R5").Value = 740.696613225348
$ws.Range("S5").Value = 0.07740146682806061
$ws.Range("T5").Value = 0.07740146682806059

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.308268
$ws.Range("H6").Value = 3.924804
$ws.Range("I6").Value = 0.08239777620284613
$ws.Range("J6").Value = 0.08239777620284613
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.253965666666666
$ws.Range("N6").Value = 6.761896999999999
$ws.Range("O6").Value = 0.03365734588711396
$ws.Range("P6").Value = 0.03365734588711396
$ws.Range("Q6").Value = 2.948791154798666
$ws.Range("R6").Value = 26.539120393188
$ws.Range("S6").Value = 0.0027732904539882
$ws.Range("T6").Value = 0.0027732904539882

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.308268
$ws.Range("H7").Value = 3.924804
$ws.Range("I7").Value = 0.08239777620284613
$ws.Range("J7").Value = 0.08239777620284613
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.806737666666667
$ws.Range("N7").Value = 5.420213
$ws.Range("O7").Value = 0.02697911306883729
$ws.Range("P7").Value = 0.02697911306883729
$ws.Range("Q7").Value = 2.363697073694667
$ws.Range("R7").Value = 21.273273663252
$ws.Range("S7").Value = 0.002223018920797336
$ws.Range("T7").Value = 0.002223018920797336

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.499176666666667
$ws.Range("H8").Value = 7.49753
$ws.Range("I8").Value = 0.1574039873109905
$ws.Range("J8").Value = 0.1574039873109906
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 62.90731233333333
$ws.Range("N8").Value = 188.721937
$ws.Range("O8").Value = 0.9393635410440488
$ws.Range("P8").Value = 0.9393635410440487
$ws.Range("Q8").Value = 157.2164871461789
$ws.Range("R8").Value = 1414.94838431561
$ws.Range("S8").Value = 0.1478595668949046
$ws.Range("T8").Value = 0.1478595668949046

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.499176666666667
$ws.Range("H9").Value = 7.49753
$ws.Range("I9").Value = 0.1574039873109905
$ws.Range("J9").Value = 0.1574039873109906
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.253965666666666
$ws.Range("N9").Value = 6.761896999999999
$ws.Range("O9").Value = 0.03365734588711396
$ws.Range("P9").Value = 0.03365734588711396
$ws.Range("Q9").Value = 5.63305840160111
$ws.Range("R9").Value = 50.69752561441
$ws.Range("S9").Value = 0.005297800444936905
$ws.Range("T9").Value = 0.005297800444936906

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.499176666666667
$ws.Range("H10").Value = 7.49753
$ws.Range("I10").Value = 0.1574039873109905
$ws.Range("J10").Value = 0.1574039873109906
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.806737666666667
$ws.Range("N10").Value = 5.420213
$ws.Range("O10").Value = 0.02697911306883729
$ws.Range("P10").Value = 0.02697911306883729
$ws.Range("Q10").Value = 4.515356619321111
$ws.Range("R10").Value = 40.63820957389
$ws.Range("S10").Value = 0.004246619971149043
$ws.Range("T10").Value = 0.004246619971149044
